# Commit: "Test data for Austria market added"
#
# Adds three new market test-data tabs - Netherlands, Austria and Denmark -
# built from the existing "Greece" tab template (same layout every other
# market sheet in this workbook follows).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Netherlands: straight clone of Greece's layout.
# ---------------------------------------------------------------------
$greece = $wb.Worksheets.Item("Greece")
$greece.Copy([System.Reflection.Missing]::Value, $greece)
$nl = $wb.Worksheets.Item($wb.Worksheets.Count)
$nl.Name = "Netherlands"

# Narrower data columns on the new market sheets.
$nl.Columns("B").ColumnWidth = 15.21875
$nl.Columns("C").ColumnWidth = 12.21875
$nl.Columns("D").ColumnWidth = 13.6640625

# Product code set before the market name so new shared strings land in the
# same order as the source edit.
$nl.Range("B4").Value = "NGC-3144/T2199"
$nl.Range("B2").Value = "Netherlands Market"
$nl.Rows("2").RowHeight = 28.8

$nl.Range("H22").Select()

# ---------------------------------------------------------------------
# Austria: clone of Netherlands, minus the "PR1D2-Unmonitored" row.
# ---------------------------------------------------------------------
$nl.Copy([System.Reflection.Missing]::Value, $nl)
$at = $wb.Worksheets.Item($wb.Worksheets.Count)
$at.Name = "Austria"

$at.Rows("10").Delete()

$at.Range("B4").Value = "NGC-3817/T2306"
$at.Range("B2").Value = "Austria Market"

$at.Range("H22").Select()

# ---------------------------------------------------------------------
# Denmark: clone of Netherlands again (keeps all six accessory rows), with
# "MZX Bezel Large" and "PR1D2-Unmonitored" swapped.
# ---------------------------------------------------------------------
$nl.Copy([System.Reflection.Missing]::Value, $at)
$dk = $wb.Worksheets.Item($wb.Worksheets.Count)
$dk.Name = "Denmark"

$dk.Range("A10").Value = "MZX Bezel Large"
$dk.Range("A11").Value = "PR1D2-Unmonitored"

$dk.Range("B4").Value = "NGC-2913/T2798"
$dk.Range("B2").Value = "Denmark Market"

$dk.Range("E19").Select()

# Denmark is the last-edited / active tab.
$dk.Activate()
